$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting the existing rows (old row 3 ->
# row 4, old row 4 -> row 5, etc.) down by one.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the new use case.
$ws.Range("B3").Value = "See a List of Tasks By Status"
$ws.Range("C3").Value = "GET"

# Match the saved selection state (active cell B4).
$ws.Range("B4").Select()
